$d = $word.ActiveDocument

# 1. Title change (appears twice: Heading1 at top, and bold line near bottom)
$d.Content.Find.Execute(
    "Play Big Thunder King Strike for Free - Ainsworth Slot Game",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Play Big Thunder King Strike Free: Review & Gameplay", 2) | Out-Null

# 2. Insert a new bullet "Variety of symbols and features" right before the
#    "Free spins round with multipliers" bullet in the "What we like" list.
$targetIndex = 0
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -eq "Free spins round with multipliers`r") {
        $targetIndex = $i
        break
    }
}
if ($targetIndex -gt 0) {
    $p = $d.Paragraphs.Item($targetIndex)
    $rng = $p.Range
    $rng.InsertParagraphBefore()
    $newp = $d.Paragraphs.Item($targetIndex)
    $sel = $word.Selection
    $sel.SetRange($newp.Range.Start, $newp.Range.Start)
    $sel.TypeText("Variety of symbols and features")
}

# 3. Replace "Generous with free spins" bullet text with "Potential for big
#    winnings" (the following bullet, "Several ways to win through
#    combinations, multipliers, and jackpots", is removed entirely).
$targetIndex = 0
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -eq "Generous with free spins`r") {
        $targetIndex = $i
        break
    }
}
if ($targetIndex -gt 0) {
    $p = $d.Paragraphs.Item($targetIndex)
    $nextp = $d.Paragraphs.Item($targetIndex + 1)
    if ($nextp.Range.Text -eq "Several ways to win through combinations, multipliers, and jackpots`r") {
        $fullRange = $d.Range($p.Range.Start, $nextp.Range.End)
        $fullRange.Text = "Potential for big winnings`r"
    }
}

# 4. "Slightly below average RTP" -> "Slightly below-average RTP"
$d.Content.Find.Execute(
    "Slightly below average RTP",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Slightly below-average RTP", 2) | Out-Null

# 5. "Limited availability of gorilla Wild symbols" -> "Limited availability
#    of Wild symbol"
$d.Content.Find.Execute(
    "Limited availability of gorilla Wild symbols",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Limited availability of Wild symbol", 2) | Out-Null

# 6. Italic summary paragraph near the end.
$d.Content.Find.Execute(
    "Explore the immersive jungle theme of Big Thunder King Strike by Ainsworth with free spins, multipliers, and jackpots. Play for free and win big!",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Play Big Thunder King Strike for free and enjoy immersive jungle-themed gameplay with potential for big winnings.", 2) | Out-Null
